$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 9 new rows starting at row 5 ---
# This pushes the old row 5 ("ACTest" data row) -> row 14
# and the old row 6 ("Final Total" row) -> row 15.
$ws.Range("A5:A13").EntireRow.Insert()

# --- 2. Row 4 becomes the "Compilation" row: drop its old B/C header cells ---
$ws.Range("B4:C4").Clear()
$ws.Range("A4").Value2 = "Compilation"

# --- 3. Remove stray ghost cells the row-insert left behind in column C/D ---
$ws.Range("C5:D13").Clear()

# --- 4. Fill the new compilation rows (5-12): ClassName | " Complies" ---
# Use the existing non-bold data-row style (row 2) as the format source so we
# reuse the workbook's existing cellXf (s="2") instead of creating new ones.
$ws.Range("A2").Copy()
$ws.Range("A5:B12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$classes = @("Device.java", "PortableDevice.java", "AC.java", "Fan.java", "StandingFan.java", "CeilingFan.java", "Room.java", "CoolingSimulation.java")
for ($i = 0; $i -lt $classes.Length; $i++) {
    $r = 5 + $i
    $ws.Cells.Item($r, 1).Value2 = $classes[$i]
    $ws.Cells.Item($r, 2).Value2 = " Complies"
}

# --- 5. Row 13: fresh "Test Class" / "Test Method" / "Comment" header ---
# Use the bold header style (row 1, which carries s="1") is too large (14pt);
# row 4 (A4) still carries the bold 12pt style (s="3") used by this header.
$ws.Range("A4").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A13").Value2 = "Test Class"
$ws.Range("B13").Value2 = "Test Method"
$ws.Range("C13").Value2 = "Comment"

# --- 6. Column A widens to fit the longer class names ---
$ws.Columns.Item(1).ColumnWidth = 22.43

$excel.CutCopyMode = 0
